# Auto-generated edit script applying Universalis price-refresh values
# to the Seraph_Profits Leve-profit tables across all class sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 250
$ws.Range("J9").Value = 250
$ws.Range("L9").Value = 250
$ws.Range("N9").Value = -588
$ws.Range("H18").Value = 10876.829
$ws.Range("I18").Value = 10918.966
$ws.Range("J18").Value = 10775
$ws.Range("K18").Value = 10918.966
$ws.Range("L18").Value = 10775
$ws.Range("M18").Value = -10634.966
$ws.Range("N18").Value = -11343
$ws.Range("H19").Value = 2366.3333
$ws.Range("J19").Value = 2657.9167
$ws.Range("L19").Value = 2657.9167
$ws.Range("N19").Value = -3007.9167
$ws.Range("H132").Value = 2787.0588
$ws.Range("I132").Value = 1327.9286
$ws.Range("K132").Value = 3983.7858
$ws.Range("M132").Value = -1453.7858

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").ClearContents()
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = 0
$ws.Range("H61").Value = 6448.091
$ws.Range("J61").Value = 3505.5
$ws.Range("L61").Value = 3505.5
$ws.Range("N61").Value = -3929.5
$ws.Range("H102").Value = 13890565
$ws.Range("I102").Value = 15874860
$ws.Range("K102").Value = 15874860
$ws.Range("M102").Value = -15873238
$ws.Range("H136").Value = 6448.091
$ws.Range("J136").Value = 3505.5
$ws.Range("L136").Value = 10516.5
$ws.Range("N136").Value = -15616.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").ClearContents()
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = 0
$ws.Range("I94").Value = 1683.3334
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 1683.3334
$ws.Range("L94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -1232.3334
$ws.Range("H99").Value = 1570
$ws.Range("I99").Value = 1618.3846
$ws.Range("K99").Value = 1618.3846
$ws.Range("M99").Value = -120.3846000000001
$ws.Range("H105").Value = 7599.778
$ws.Range("I105").Value = 6049.875
$ws.Range("K105").Value = 6049.875
$ws.Range("M105").Value = -4302.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 890
$ws.Range("I21").Value = 890
$ws.Range("K21").Value = 890
$ws.Range("M21").Value = -655
$ws.Range("H62").Value = 268332.66
$ws.Range("I62").Value = 202499.5
$ws.Range("K62").Value = 202499.5
$ws.Range("M62").Value = -201875.5
$ws.Range("H65").Value = 268332.66
$ws.Range("I65").Value = 202499.5
$ws.Range("K65").Value = 1012497.5
$ws.Range("M65").Value = -1009377.5
$ws.Range("H99").Value = 16130.952
$ws.Range("I99").Value = 10972.875
$ws.Range("K99").Value = 10972.875
$ws.Range("M99").Value = -9474.875
$ws.Range("H107").Value = 1047.6857
$ws.Range("I107").Value = 737.1667
$ws.Range("J107").Value = 1376.4706
$ws.Range("K107").Value = 737.1667
$ws.Range("L107").Value = 1376.4706
$ws.Range("M107").Value = 1182.8333
$ws.Range("N107").Value = -5216.470600000001
$ws.Range("H126").Value = 16130.952
$ws.Range("I126").Value = 10972.875
$ws.Range("K126").Value = 32918.625
$ws.Range("M126").Value = -30448.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 2749
$ws.Range("J52").Value = 2749
$ws.Range("L52").Value = 8247
$ws.Range("N52").Value = -8779
$ws.Range("H64").Value = 3269.5
$ws.Range("I64").Value = 3269.5
$ws.Range("K64").Value = 9808.5
$ws.Range("M64").Value = -9538.5
$ws.Range("H67").Value = 3269.5
$ws.Range("I67").Value = 3269.5
$ws.Range("K67").Value = 9808.5
$ws.Range("M67").Value = -8872.5
$ws.Range("H98").Value = 2739.6
$ws.Range("I98").Value = 2724.5
$ws.Range("K98").Value = 8173.5
$ws.Range("M98").Value = -6675.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 15000400
$ws.Range("I11").Value = 18750000
$ws.Range("J11").Value = 2000
$ws.Range("K11").Value = 18750000
$ws.Range("L11").Value = 2000
$ws.Range("M11").Value = -18749861
$ws.Range("N11").Value = -2278
$ws.Range("H80").Value = 3019.1924
$ws.Range("I80").Value = 2979.125
$ws.Range("J80").Value = 3500
$ws.Range("K80").Value = 2979.125
$ws.Range("L80").Value = 3500
$ws.Range("M80").Value = -1981.125
$ws.Range("N80").Value = -5496
$ws.Range("H83").Value = 3019.1924
$ws.Range("I83").Value = 2979.125
$ws.Range("J83").Value = 3500
$ws.Range("K83").Value = 14895.625
$ws.Range("L83").Value = 17500
$ws.Range("M83").Value = -9903.625
$ws.Range("N83").Value = -27484
$ws.Range("H122").Value = 45454.24
$ws.Range("I122").Value = 4328.25
$ws.Range("J122").Value = 118567.11
$ws.Range("K122").Value = 12984.75
$ws.Range("L122").Value = 355701.33
$ws.Range("M122").Value = -10534.75
$ws.Range("N122").Value = -360601.33
$ws.Range("H126").Value = 4506.4116
$ws.Range("I126").Value = 3687.2856
$ws.Range("K126").Value = 11061.8568
$ws.Range("M126").Value = -8591.856800000001
$ws.Range("H132").Value = 3285.5715
$ws.Range("J132").Value = 3399.8
$ws.Range("L132").Value = 10199.4
$ws.Range("N132").Value = -15259.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1173.8
$ws.Range("J22").Value = 1450
$ws.Range("L22").Value = 1450
$ws.Range("N22").Value = -2040
$ws.Range("H27").Value = 1173.8
$ws.Range("J27").Value = 1450
$ws.Range("L27").Value = 1450
$ws.Range("N27").Value = -1664
$ws.Range("H40").Value = 3149.3333
$ws.Range("I40").Value = 3149.3333
$ws.Range("K40").Value = 3149.3333
$ws.Range("M40").Value = -3013.3333
$ws.Range("H46").Value = 2800
$ws.Range("J46").Value = 2500
$ws.Range("L46").Value = 2500
$ws.Range("N46").Value = -2876
$ws.Range("H55").Value = 621.7646999999999
$ws.Range("I55").Value = 594.9
$ws.Range("K55").Value = 594.9
$ws.Range("M55").Value = -421.9
$ws.Range("H122").Value = 4503.1
$ws.Range("I122").Value = 2609.7
$ws.Range("K122").Value = 7829.099999999999
$ws.Range("M122").Value = -5379.099999999999
$ws.Range("H124").Value = 63429
$ws.Range("J124").Value = 63429
$ws.Range("L124").Value = 63429
$ws.Range("N124").Value = -73249

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1030
$ws.Range("J113").Value = 2750
$ws.Range("L113").Value = 8250
$ws.Range("N113").Value = -12590
